$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 33 - Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 659.5143
$ws.Range("I33").Value = 542.76666
$ws.Range("J33").Value = 1360
$ws.Range("K33").Value = 542.76666
$ws.Range("L33").Value = 1360
$ws.Range("M33").Value = -313.76666
$ws.Range("N33").Value = -1818

# Row 49 - Going Nowhere Fast / Paralyzing Potion
$ws.Range("H49").Value = 1750
$ws.Range("I49").Value = 3000
$ws.Range("J49").Value = 500
$ws.Range("K49").Value = 9000
$ws.Range("L49").Value = 1500
$ws.Range("M49").Value = -8864
$ws.Range("N49").Value = -1772

# Row 132
$ws.Range("H132").Value = 4547.0625
$ws.Range("I132").Value = 2765.8696
$ws.Range("K132").Value = 8297.6088
$ws.Range("M132").Value = -5767.6088

# Row 138 - All-night Crafting / Cunning Craftsman's Tisane (M138 removed)
$ws.Range("H138").Value = 2981.25
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2981.25
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 8943.75
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -19223.75

# ---------------------------------------------------------------------------
# ARM sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 2108.38
$ws.Range("I32").Value = 1638.9678
$ws.Range("J32").Value = 8344.857
$ws.Range("K32").Value = 1638.9678
$ws.Range("L32").Value = 8344.857
$ws.Range("M32").Value = -1351.9678
$ws.Range("N32").Value = -8918.857

# Row 61
$ws.Range("H61").Value = 240817.36
$ws.Range("I61").Value = 201478.77
$ws.Range("J61").Value = 306381.66
$ws.Range("K61").Value = 201478.77
$ws.Range("L61").Value = 306381.66
$ws.Range("M61").Value = -201266.77
$ws.Range("N61").Value = -306805.66

# Row 136
$ws.Range("H136").Value = 240817.36
$ws.Range("I136").Value = 201478.77
$ws.Range("J136").Value = 306381.66
$ws.Range("K136").Value = 604436.3099999999
$ws.Range("L136").Value = 919144.98
$ws.Range("M136").Value = -601886.3099999999
$ws.Range("N136").Value = -924244.98

# ---------------------------------------------------------------------------
# BSM sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 99
$ws.Range("H99").Value = 7695962.5
$ws.Range("I99").Value = 3501943.2
$ws.Range("J99").Value = 14286564
$ws.Range("K99").Value = 3501943.2
$ws.Range("L99").Value = 14286564
$ws.Range("M99").Value = -3500445.2
$ws.Range("N99").Value = -14289560

# Row 107
$ws.Range("H107").Value = 1706.1086
$ws.Range("I107").Value = 1411.258
$ws.Range("J107").Value = 2315.4666
$ws.Range("K107").Value = 1411.258
$ws.Range("L107").Value = 2315.4666
$ws.Range("M107").Value = 508.742
$ws.Range("N107").Value = -6155.4666

# ---------------------------------------------------------------------------
# CRP sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 3406.375
$ws.Range("I31").Value = 2464.2856
$ws.Range("K31").Value = 2464.2856
$ws.Range("M31").Value = -2169.2856

# Row 34
$ws.Range("H34").Value = 3406.375
$ws.Range("I34").Value = 2464.2856
$ws.Range("K34").Value = 2464.2856
$ws.Range("M34").Value = -2262.2856

# Row 58
$ws.Range("H58").Value = 3186.6667
$ws.Range("I58").Value = 3663.9429
$ws.Range("J58").Value = 2307.4736
$ws.Range("K58").Value = 3663.9429
$ws.Range("L58").Value = 2307.4736
$ws.Range("M58").Value = -3460.9429
$ws.Range("N58").Value = -2713.4736

# Row 132
$ws.Range("H132").Value = 1952.3024
$ws.Range("I132").Value = 1188.3
$ws.Range("J132").Value = 3715.3845
$ws.Range("K132").Value = 3564.9
$ws.Range("L132").Value = 11146.1535
$ws.Range("M132").Value = -1034.9
$ws.Range("N132").Value = -16206.1535

# Row 134
$ws.Range("H134").Value = 1676.1111
$ws.Range("I134").Value = 1145.0488
$ws.Range("J134").Value = 2378.484
$ws.Range("K134").Value = 3435.1464
$ws.Range("L134").Value = 7135.451999999999
$ws.Range("M134").Value = -900.1464000000001
$ws.Range("N134").Value = -12205.452

# Row 136
$ws.Range("H136").Value = 3186.6667
$ws.Range("I136").Value = 3663.9429
$ws.Range("J136").Value = 2307.4736
$ws.Range("K136").Value = 10991.8287
$ws.Range("L136").Value = 6922.4208
$ws.Range("M136").Value = -8441.8287
$ws.Range("N136").Value = -12022.4208

# ---------------------------------------------------------------------------
# CUL sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 126
$ws.Range("H126").Value = 2133.3333
$ws.Range("I126").Value = 466.66666
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 1399.99998
$ws.Range("L126").Value = 11400
$ws.Range("M126").Value = 3540.00002
$ws.Range("N126").Value = -21280

# Row 129
$ws.Range("H129").Value = 1556.0625
$ws.Range("I129").Value = 549.25
$ws.Range("J129").Value = 2562.875
$ws.Range("K129").Value = 1647.75
$ws.Range("L129").Value = 7688.625
$ws.Range("M129").Value = 3352.25
$ws.Range("N129").Value = -17688.625

# Row 130
$ws.Range("H130").Value = 500481.5
$ws.Range("I130").Value = 930
$ws.Range("K130").Value = 2790
$ws.Range("M130").Value = 2230

# Row 132
$ws.Range("H132").Value = 8427.182000000001
$ws.Range("I132").Value = 5300
$ws.Range("J132").Value = 11033.167
$ws.Range("K132").Value = 47700
$ws.Range("L132").Value = 99298.503
$ws.Range("M132").Value = -45170
$ws.Range("N132").Value = -104358.503

# Row 137
$ws.Range("H137").Value = 1741.25
$ws.Range("I137").Value = 1847.1428
$ws.Range("K137").Value = 5541.428400000001
$ws.Range("M137").Value = -441.4284000000007

# Row 141
$ws.Range("H141").Value = 4493.095
$ws.Range("I141").Value = 5285.5
$ws.Range("J141").Value = 3772.7273
$ws.Range("K141").Value = 15856.5
$ws.Range("L141").Value = 11318.1819
$ws.Range("M141").Value = -10676.5
$ws.Range("N141").Value = -21678.1819

# ---------------------------------------------------------------------------
# GSM sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 5561.9707
$ws.Range("I70").Value = 3895.689
$ws.Range("J70").Value = 8822.087
$ws.Range("K70").Value = 3895.689
$ws.Range("L70").Value = 8822.087
$ws.Range("M70").Value = -3625.689
$ws.Range("N70").Value = -9362.087

# Row 73
$ws.Range("H73").Value = 5561.9707
$ws.Range("I73").Value = 3895.689
$ws.Range("J73").Value = 8822.087
$ws.Range("K73").Value = 3895.689
$ws.Range("L73").Value = 8822.087
$ws.Range("M73").Value = -2959.689
$ws.Range("N73").Value = -10694.087

# ---------------------------------------------------------------------------
# WVR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 96
$ws.Range("H96").Value = 2340.0625
$ws.Range("I96").Value = 748.3333
$ws.Range("J96").Value = 3295.1
$ws.Range("K96").Value = 748.3333
$ws.Range("L96").Value = 3295.1
$ws.Range("M96").Value = 624.6667
$ws.Range("N96").Value = -6041.1

# Row 132
$ws.Range("H132").Value = 1691.4315
$ws.Range("I132").Value = 1159.2554
$ws.Range("J132").Value = 2212.5208
$ws.Range("K132").Value = 3477.7662
$ws.Range("L132").Value = 6637.562399999999
$ws.Range("M132").Value = -947.7662
$ws.Range("N132").Value = -11697.5624
